# Auto-generated Word COM-interop script
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false

function ReplaceText($d, $old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "FAILED REPLACE: $old"
        return $false
    }
    $r.Text = $new
    return $true
}

function InsertAfterText($d, $anchor, $newText) {
    $r = $d.Content
    $ok = $r.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "FAILED FIND (insert after): $anchor"
        return
    }
    $ins = $d.Range($r.End, $r.End)
    $ins.InsertAfter($newText)
}

function DeleteText($d, $old) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "FAILED DELETE: $old"
        return
    }
    $r.Delete()
}

$d = $word.ActiveDocument

# --- Global font fix: TimesNewToman -> Times New Roman ---
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Font.Name = "TimesNewToman"
$find.Replacement.Font.Name = "Times New Roman"
$find.Execute("", $false, $false, $false, $false, $false, $true, 1, $true, "", 2) | Out-Null

# Title
ReplaceText $d 'The Profound Enigma of Time''s Flow' 'The Convergence of Creativity and Reality: An Interwoven Journey Through Art History' | Out-Null

# Author first word
ReplaceText $d 'Dr' 'Prof' | Out-Null

# Author full name
ReplaceText $d ' Alex Richards' ' Eleanor Hayes' | Out-Null

# Email local part
ReplaceText $d 'richards' 'eleanorhayes@gmail' | Out-Null

# Email domain part
ReplaceText $d 'alex@researchhub' 'com' | Out-Null

# Remove trailing .edu
DeleteText $d '.edu'

# Body sentence 1 -> '1'
ReplaceText $d 'In the boundless tapestry of human knowledge, few concepts inspire as much fascination and perplexity as time' '1' | Out-Null

# Body sentence 2
ReplaceText $d ' Its relentless progression, its unfathomable nature, and its profound implications for our perception of reality have captivated thinkers for eons' ' From the cave walls of Lascaux to the vibrant streets of modern day, art has been humanity''s unwavering companion, mirroring our beliefs, triumphs, woes, and dreams' | Out-Null

# Body sentence 3
ReplaceText $d ' Time''s arrow, with its inexorable forward motion, marking the passage from past to present to future, remains an enigma that defies complete comprehension' ' Painted ceilings and sculpted verses narrate sagas of our evolution, holding up a looking glass to the ever-changing visage of the world we build' | Out-Null

# Insert new sentence before first <br/>
InsertAfterText $d ' Painted ceilings and sculpted verses narrate sagas of our evolution, holding up a looking glass to the ever-changing visage of the world we build.' ' A testament to human ingenuity and resilience, art weaves the intricate tapestry of our cultural heritage, earning its place as an integral part of our lives.'

# Body sentence after br -> '2'
ReplaceText $d 'As we traverse the corridors of time, its passing leaves an imprint on our consciousness, etching memories, shaping our understanding of existence, and dictating the rhythms of our lives' '2' | Out-Null

# Body sentence 5
ReplaceText $d ' The ephemerality of time imbues each moment with a heightened sense of significance, compelling us to ponder our own mortality and the fleeting nature of our experiences' ' Art transcends time, offering a medium that bridges the gap between disparate generations' | Out-Null

# Body sentence 6
ReplaceText $d ' Yet, despite its omnipresence, time eludes easy definition, slipping through our grasp like grains of sand' ' The strokes of a maestro painter, crafted centuries ago, continue to hold audiences spellbound, stirring emotions and shedding light onto our shared human experience' | Out-Null

# Insert two new sentences before second <br/>
InsertAfterText $d ' The strokes of a maestro painter, crafted centuries ago, continue to hold audiences spellbound, stirring emotions and shedding light onto our shared human experience.' ' A symphony''s haunting melodies, composed decades past, echo through concert halls, enchanting listeners with their timeless beauty. Art''s kaleidoscopic expressions provide a portal, allowing us to connect with the hopes, fears and dreams of those who came before us.'

# Body sentence after 2nd br -> '3'
ReplaceText $d 'The relentless march of time has been a muse for poets, philosophers, and scientists alike, each attempting to unravel its mysteries' '3' | Out-Null

# Body sentence 9 (last original)
ReplaceText $d ' From the poetic musings of William Wordsworth, who captured the fleeting beauty of time in his "Ode: Intimations of Immortality," to the philosophical ponderings of Henri Bergson, who explored the nature of time and duration in his seminal work "Time and Free Will," humanity''s quest to understand time has been a relentless pursuit' ' Art provokes reflection and introspection, pushing boundaries and challenging established conventions' | Out-Null

# Insert two new trailing sentences at end of paragraph
InsertAfterText $d ' Art provokes reflection and introspection, pushing boundaries and challenging established conventions.' ' It holds up a mirror to our social, cultural and political realities, compelling us to confront uncomfortable truths and confront issues that plague our societies. Art can ignite change, amplify voices, and foster empathy, challenging us to rethink and reconstruct both ourselves and the world around us.'

# Summary sentence 1
ReplaceText $d 'Time remains an enigma that captivates and confounds, its nature defying simple explanation' 'Art, in its dazzling myriad forms, serves as a potent chronicle of human civilization' | Out-Null

# Summary sentence 2
ReplaceText $d ' The relentless progression of time, its role in shaping consciousness and perception, and its profound implications for human existence have inspired countless works of art, philosophy, and scientific exploration' ' Beyond its aesthetic allure, it offers a profound means of comprehending our past, making sense of the present, and envisioning a better future' | Out-Null

# Summary sentence 3
ReplaceText $d ' As we continue to probe the depths of this timeless mystery, we embrace the awe-inspiring beauty of its elusiveness, acknowledging that time''s profound enigma may forever surpass our grasp' ' Art invites us on a journey of discovery, both inward and outward, bridging the gap between generations and cultures' | Out-Null

# Insert new trailing summary sentence
InsertAfterText $d ' Art invites us on a journey of discovery, both inward and outward, bridging the gap between generations and cultures.' ' It challenges societal norms, ignites change, and remains a pivotal force in shaping the ever-evolving tapestry of human existence.'

# --- Add trailing empty paragraph at the very end of the document ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter() | Out-Null

Write-Output "DONE"